# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the existing header formatting (bold, centered, thin border) instead
# of minting a new style - copy the format from the last existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2 through 56) ---------------------------------------------
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 60   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 102  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
